$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos price/volume refresh (GitHub Actions scheduled update).
# Column D occasionally holds numeric-looking text (e.g. "1.00", "54.30") that
# must stay literal text (matches source feed formatting) instead of being
# auto-coerced to a number by Excel, so those cells are explicitly forced to
# the Text number format before the value is written.

# Row 2
$ws.Range('D2').Value = '66.497.12'
$ws.Range('E2').Value = '  +8.32%  '

# Row 3
$ws.Range('D3').Value = '3.466.53'
$ws.Range('E3').Value = '  +12.04%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.28%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.18'
$ws.Range('E5').Value = '  +12.56%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '544.59'
$ws.Range('E6').Value = '  +7.55%  '

# Row 7
$ws.Range('D7').Value = '3.456.31'
$ws.Range('E7').Value = '  +11.65%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.602'
$ws.Range('E8').Value = '  +4.41%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.08%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.626'
$ws.Range('E10').Value = '  +8.49%  '

# Row 11
$ws.Range('E11').Value = '  +19.21%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.30'
$ws.Range('E12').Value = '  +6.17%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000265'
$ws.Range('E13').Value = '  +10.53%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.28'
$ws.Range('E14').Value = '  +7.38%  '

# Row 15
$ws.Range('D15').Value = '4.023.58'
$ws.Range('E15').Value = '  +12.63%  '

# Row 16
$ws.Range('D16').Value = '3.471.78'
$ws.Range('E16').Value = '  +12.66%  '

# Row 17
$ws.Range('E17').Value = '  +8.06%  '

# Row 18
$ws.Range('D18').Value = '66.606.45'
$ws.Range('E18').Value = '  +9.15%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.97'
$ws.Range('E19').Value = '  +8.83%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('E20').Value = '  +10.87%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.985'
$ws.Range('E21').Value = '  +6.71%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '417.10'
$ws.Range('E22').Value = '  +18.18%  '

# Row 23
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.20'
$ws.Range('E23').Value = '  +11.56%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.00'
$ws.Range('E24').Value = '  +7.60%  '

# Row 25
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.86'
$ws.Range('E25').Value = '  +7.69%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.09'
$ws.Range('E26').Value = '  +4.20%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.87'
$ws.Range('E27').Value = '  +14.28%  '

# Row 28
$ws.Range('E28').Value = '  +0.93%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.78'
$ws.Range('E29').Value = '  +9.88%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.71'
$ws.Range('E30').Value = '  +11.91%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.82'
$ws.Range('E31').Value = '  +9.77%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '653.48'
$ws.Range('E32').Value = '  +4.84%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.58'
$ws.Range('E33').Value = '  +7.25%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.62'
$ws.Range('E34').Value = '  +7.21%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.109'
$ws.Range('E35').Value = '  +9.15%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.90'
$ws.Range('E36').Value = '  +5.94%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.05'
$ws.Range('E37').Value = '  +9.08%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0808'
$ws.Range('E38').Value = '  +22.59%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.29%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.386'
$ws.Range('E40').Value = '  +7.12%  '

# Row 41
$ws.Range('E41').Value = '  +15.46%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.32'
$ws.Range('E42').Value = '  +17.96%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.58%  '

# Row 44
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.45'
$ws.Range('E44').Value = '  +21.25%  '

# Row 45
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.990.09'
$ws.Range('E45').Value = '  +8.38%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.61'
$ws.Range('E46').Value = '  +8.02%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.87'
$ws.Range('E47').Value = '  +17.59%  '

# Row 48
$ws.Range('E48').Value = '  +10.51%  '

# Row 49
$ws.Range('E49').Value = '  +4.39%  '

# Row 50
$ws.Range('E50').Value = '  +20.21%  '

# Row 51
$ws.Range('E51').Value = '  +8.60%  '
